$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) tweaks
$ws.Range("B2").Value = 35.761089351678933
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 39.866832231383704
$ws.Range("E2").Value = 44.53980716691045

# Row 3 (STR) tweaks
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 50.302108950983936
$ws.Range("D3").Value = 46.752859715113551
$ws.Range("E3").Value = 44.76646235616586

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
